# إضافة حدث جديد في Card22 by admin at 2025-12-08 08:08:17
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

$nanCols = @("B","C","D","E","F","G","H","I","J","K")

# Row 19: columns B..K were blank inline-string cells; they become the text "nan"
foreach ($col in $nanCols) {
    $ws.Range("$col`19").Value = "nan"
}

# Row 20: brand-new service event row
# Column A holds the card number as text ("22"); a leading apostrophe forces
# Excel to store the numeric-looking value as text instead of a number.
$ws.Range("A20").Value = "'22"

# Columns B..K stay blank, just like on every other row of this sheet.

$ws.Range("L20").Value = "1\9\2025"
$ws.Range("M20").Value = "800 t"
$ws.Range("N20").Value = "تم تغيير الجرائد الخلفيه (1_5_8)ومعايره"
$ws.Range("O20").Value = "الخبير"
